# Prepend a new "today" row of price data to the top of the data table
# (row 2, just below the header row) and push all existing rows down by
# one, matching the source site's daily scrape.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows (2:13) down to make room for the new entry.
$ws.Rows("2:2").Insert()

# The inserted row inherits the header row's bold/centered/bordered
# formatting from the row above; strip that back to the plain "no
# explicit style" look used by the rest of the data rows.
$ws.Range("A2:D2").ClearFormats()

# Force column A to be treated as text so the date string isn't
# auto-converted into a date serial number (matches the other rows,
# which all store the date as plain text).
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2025-12-03"
$ws.Range("B2").Value = 783.5
$ws.Range("C2").Value = 1112
$ws.Range("D2").Value = 3610

# Drop the temporary text number-format override so the new row ends up
# with no explicit style, just like the rest of the data rows.
$ws.Range("A2").ClearFormats()
